# Apply the edit described by the commit "bug fixi rename some indicators":
#  - rows 2-6 (D column): Thermal Output Capacity default values -> 0.001 (1E-3)
#  - rows 20-25 and 31-35 (G:H columns): populate min/max with +/- huge numbers,
#    using the same scientific-notation style already used elsewhere in the sheet
#  - rows 46-50 (F column): rename unit from "EUR/MW" to "EUR/(MW*yr)"
#  - move the active selection to L44

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- D2:D6 -> 0.001 -------------------------------------------------------
$ws.Range("D2").Value = 0.001
$ws.Range("D3").Value = 0.001
$ws.Range("D4").Value = 0.001
$ws.Range("D5").Value = 0.001
$ws.Range("D6").Value = 0.001

# --- G/H min-max for rows 20-25 and 31-35 ---------------------------------
$bigNumber = 1 * [math]::Pow(10, 102)
$minMaxRows = 20,21,22,23,24,25,31,32,33,34,35
foreach ($r in $minMaxRows) {
    $ws.Cells.Item($r, 7).Value = -$bigNumber
    $ws.Cells.Item($r, 8).Value = $bigNumber
}

# Re-apply the scientific number format (with its quote-prefix flag) that is
# already used by the sibling cells (e.g. H2) elsewhere in the sheet, so the
# new cells end up sharing the same cell style as the rest of the min/max
# columns instead of a plain General style.
$ws.Range("H2").Copy()
foreach ($r in $minMaxRows) {
    $target = $ws.Range($ws.Cells.Item($r, 7), $ws.Cells.Item($r, 8))
    $target.PasteSpecial(-4122)
}
$excel.CutCopyMode = $false

# --- F46:F50 unit rename ---------------------------------------------------
$ws.Range("F46").Value = "EUR/(MW*yr)"
$ws.Range("F47").Value = "EUR/(MW*yr)"
$ws.Range("F48").Value = "EUR/(MW*yr)"
$ws.Range("F49").Value = "EUR/(MW*yr)"
$ws.Range("F50").Value = "EUR/(MW*yr)"

# --- move selection ----------------------------------------------------
[void]$ws.Range("L44").Select()
